# Action_overview_catheterisation_Women.xlsx
# "lighting for catheterisation_wmoen scene"
#
# Adds new task/issue rows to the "Blad2" (issues) sheet, highlights the
# in-progress rows with Excel's built-in "Neutral" cell style, clears a
# stray yellow highlight on "Blad1", and refreshes the active selections.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Blad1
$ws2 = $wb.Worksheets.Item(2)   # Blad2

# ---------------------------------------------------------------------
# Blad1: the earlier highlighted note on E12 is no longer relevant -
# drop the yellow fill back to the sheet default.
# ---------------------------------------------------------------------
$ws1.Range("E12").ClearFormats()

# ---------------------------------------------------------------------
# Blad2: fill in the remaining cells of the in-progress issue table and
# append new rows describing follow-up work.
# ---------------------------------------------------------------------

# Row 2 (A2:E2) already holds its text - only the highlight below is new.

# Row 3: finish the issue (B3 already held placeholder text "Creating ").
$ws2.Range("B3").Value = "Creating move positions for GauzeTrayWet & PlasticTrashbucket"
$ws2.Range("C3").Value = "GauzeTrayWet & PlasticTrashbucket need to be moved within the protocol. To do this, we need multiple drop locations. Please do this for GauzeTrayWet & PlasticTrashbucket. How to do this can be scene here: https://github.com/GijsTempel/care-up/wiki/Objects under title ' moving objects'"
$ws2.Range("D3").Value = "Implement, medium priority develop"
$ws2.Range("E3").Value = "Dani"

# Row 4 (B4:E4) already has its values (Detecting move.../Alexander) - untouched.

# Row 5: new issue.
$ws2.Range("B5").Value = "Creating holding animation for PlasticTrashbucket"
$ws2.Range("C5").Value = "In this protocol we need to move the PasticTrashbucket. Please create a holding animaton for this"
$ws2.Range("D5").Value = "Animation, medium priority develop"
$ws2.Range("E5").Value = "Vitalii"

# Row 6: new issue.
$ws2.Range("B6").Value = "Adding Place Gauze tray  & plastic trashbucket. on pad step to XML"
$ws2.Range("C6").Value = "Adding Place Gauze tray  & plastic trashbucket. on pad to XML actions list. Create to seperate actions for PlasticTrashBucket and GauzeTray so players can move it in custom order. "
$ws2.Range("D6").Value = "Implement, medium priority develop"
$ws2.Range("E6").Value = "Vitalii"

# Row 8: new, entered before row 7 is finished off below.
$ws2.Range("A8").Value = "Creating animation sequence for steps (11t/m17)"
$ws2.Range("B8").Value = "Creating animation sequence for cleaning genitals"

# Row 7: new issue.
$ws2.Range("B7").Value = "Implementing holding animation issue #851"
$ws2.Range("C7").Value = "Implement holding animation of PlasticTrashbucket so it can be called in XML"
$ws2.Range("D7").Value = "Implement, medium priority develop"
$ws2.Range("E7").Value = "Vitalii"

# Highlight the whole in-progress block (rows 2-7, cols B:E) using Excel's
# built-in "Neutral" cell style (orange text on a pale-yellow fill).
$ws2.Range("B2:E7").Style = "Neutral"

# Rows 11, 10, 9 entered in that order (plain, unstyled).
$ws2.Range("B11").Value = "Adding animation sequence steps for cleaning genitals to actions XML"
$ws2.Range("B10").Value = "Creating animations sequence XML for cleaning genitals animation sequence "
$ws2.Range("B9").Value = "Implementing animation sequence so it can be called in XML "

# Widen column C to fit the new, much longer descriptions.
$ws2.Columns.Item(3).ColumnWidth = 63.43

# ---------------------------------------------------------------------
# Refresh the active selections to match where editing finished.
# ---------------------------------------------------------------------
$ws1.Range("A19").Select()
$ws2.Range("B11").Select()
